$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM line item: Cylinder (row 16), with hyperlink to the Grainger product page
$ws.Range("B16").Value = "Cylinder"
$url = "http://www.grainger.com/product/MAXIM-Hydraulic-Cylinder-6FDA8?Pid=search"
$ws.Hyperlinks.Add($ws.Range("C16"), $url)
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 236.5

# Leave the current selection where the user ended up editing
$ws.Range("B18").Select() | Out-Null
